$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9455
$ws.Range("D2").Value = 8375
$ws.Range("E2").Value = 0.8857747223691169
$ws.Range("F2").Value = 0.8836252373918548
$ws.Range("G2").Value = 0.09689384242872673
$ws.Range("H2").Value = 0.08561784451789264
$ws.Range("I2").Value = 41234906.25223832
$ws.Range("J2").Value = 14421586.78406116
$ws.Range("L2").Value = 14421586.78406116
$ws.Range("M2").Value = 55656493.03629947
$ws.Range("N2").Value = 800568137.6472001
$ws.Range("O2").Value = 782868330.6432
$ws.Range("P2").Value = 0.018014190305456
$ws.Range("Q2").Value = 0.01842147168248901

$ws.Range("C3").Value = 9647
$ws.Range("D3").Value = 8557
$ws.Range("E3").Value = 0.8870115061677205
$ws.Range("F3").Value = 0.8852679495137595
$ws.Range("G3").Value = 0.1020477647857244
$ws.Range("H3").Value = 0.09033961548432071
$ws.Range("I3").Value = 48109582.45102569
$ws.Range("J3").Value = 17588976.84787663
$ws.Range("L3").Value = 17588976.84787663
$ws.Range("M3").Value = 65698559.29890232
$ws.Range("N3").Value = 838121313.564728
$ws.Range("O3").Value = 820641137.540658
$ws.Range("P3").Value = 0.02098619443653873
$ws.Range("Q3").Value = 0.0214332136707017

$ws.Range("C4").Value = 9840
$ws.Range("D4").Value = 8723
$ws.Range("E4").Value = 0.8864837398373984
$ws.Range("F4").Value = 0.8848650841955772
$ws.Range("G4").Value = 0.1062957128902854
$ws.Range("H4").Value = 0.09405736493629131
$ws.Range("I4").Value = 54616763.50503325
$ws.Range("J4").Value = 20522613.28469532
$ws.Range("L4").Value = 20522613.28469532
$ws.Range("M4").Value = 75139376.78972858
$ws.Range("N4").Value = 876530452.3962009
$ws.Range("O4").Value = 859081504.3902471
$ws.Range("P4").Value = 0.02341346296479713
$ws.Range("Q4").Value = 0.02388901772394893

$ws.Range("C5").Value = 10032
$ws.Range("D5").Value = 8909
$ws.Range("E5").Value = 0.8880582137161085
$ws.Range("F5").Value = 0.8861149791127909
$ws.Range("G5").Value = 0.1093961613107875
$ws.Range("H5").Value = 0.09693757719492795
$ws.Range("I5").Value = 60672123.38640694
$ws.Range("J5").Value = 23217186.5988863
$ws.Range("L5").Value = 23217186.5988863
$ws.Range("M5").Value = 83889309.98529324
$ws.Range("N5").Value = 914302919.6242424
$ws.Range("O5").Value = 896816814.1612692
$ws.Range("P5").Value = 0.02539331998242774
$ws.Range("Q5").Value = 0.02588843812055389

$ws.Range("C6").Value = 10235
$ws.Range("D6").Value = 9057
$ws.Range("E6").Value = 0.884904738641915
$ws.Range("F6").Value = 0.8832650672908133
$ws.Range("G6").Value = 0.1090150774465601
$ws.Range("H6").Value = 0.0962892097165491
$ws.Range("I6").Value = 63804904.40937157
$ws.Range("J6").Value = 24412498.54468551
$ws.Range("L6").Value = 24412498.54468551
$ws.Range("M6").Value = 88217402.95405708
$ws.Range("N6").Value = 955199529.3546511
$ws.Range("O6").Value = 937607703.4814696
$ws.Range("P6").Value = 0.02555748594346461
$ws.Range("Q6").Value = 0.02603700721958497
